$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "FilesTab" query (row 4 / cell B4) dropped the `File Type` and `Breed`
# columns from its RETURN clause.
$filesQuery = @'

MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE demo.sex IN ['Female']
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $filesQuery

# Refresh the active selection on the sheet to match the re-saved view state.
$ws.Activate()
$ws.Range("B4").Select()
